$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 48, shifting everything below (including the
# "Broadband" label further down at row 208) down by one row.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row with the new entry.
$ws.Range("R48").Value = "modify processed"
$ws.Range("S48").Value = "2024-09-24 22:30:14"
